$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update existing GDP per Capita values (years 1950-2010, rows 2-62) ---
$ws.Cells.Item(2, 5).Value = "'1516"
$ws.Cells.Item(3, 5).Value = "'1549"
$ws.Cells.Item(4, 5).Value = "'1581"
$ws.Cells.Item(5, 5).Value = "'1613"
$ws.Cells.Item(6, 5).Value = "'1645"
$ws.Cells.Item(7, 5).Value = "'1677"
$ws.Cells.Item(8, 5).Value = "'1709"
$ws.Cells.Item(9, 5).Value = "'1741"
$ws.Cells.Item(10, 5).Value = "'1773"
$ws.Cells.Item(11, 5).Value = "'1804"
$ws.Cells.Item(12, 5).Value = "'1793"
$ws.Cells.Item(13, 5).Value = "'1795"
$ws.Cells.Item(14, 5).Value = "'1800"
$ws.Cells.Item(15, 5).Value = "'1747"
$ws.Cells.Item(16, 5).Value = "'1779"
$ws.Cells.Item(17, 5).Value = "'1734"
$ws.Cells.Item(18, 5).Value = "'1733"
$ws.Cells.Item(19, 5).Value = "'1790"
$ws.Cells.Item(20, 5).Value = "'1871"
$ws.Cells.Item(21, 5).Value = "'1900"
$ws.Cells.Item(22, 5).Value = "'1954"
$ws.Cells.Item(23, 5).Value = "'1986"
$ws.Cells.Item(24, 5).Value = "'1916"
$ws.Cells.Item(25, 5).Value = "'1824"
$ws.Cells.Item(26, 5).Value = "'1816"
$ws.Cells.Item(27, 5).Value = "'1795"
$ws.Cells.Item(28, 5).Value = "'1694"
$ws.Cells.Item(29, 5).Value = "'1691"
$ws.Cells.Item(30, 5).Value = "'1604"
$ws.Cells.Item(31, 5).Value = "'1714"
$ws.Cells.Item(32, 5).Value = "'1680"
$ws.Cells.Item(33, 5).Value = "'1492"
$ws.Cells.Item(34, 5).Value = "'1423"
$ws.Cells.Item(35, 5).Value = "'1395"
$ws.Cells.Item(36, 5).Value = "'1305"
$ws.Cells.Item(37, 5).Value = "'1296"
$ws.Cells.Item(38, 5).Value = "'1267"
$ws.Cells.Item(39, 5).Value = "'1258"
$ws.Cells.Item(40, 5).Value = "'1240"
$ws.Cells.Item(41, 5).Value = "'1251"
$ws.Cells.Item(42, 5).Value = "'1262"
$ws.Cells.Item(43, 5).Value = "'1160.67053364406"
$ws.Cells.Item(44, 5).Value = "'1152.64606233185"
$ws.Cells.Item(45, 5).Value = "'1154.18948369278"
$ws.Cells.Item(46, 5).Value = "'1132.1717646415"
$ws.Cells.Item(47, 5).Value = "'1129.84854810737"
$ws.Cells.Item(48, 5).Value = "'1132.44372349517"
$ws.Cells.Item(49, 5).Value = "'1152.15458292731"
$ws.Cells.Item(50, 5).Value = "'1174.22722748086"
$ws.Cells.Item(51, 5).Value = "'1206.89376146016"
$ws.Cells.Item(52, 5).Value = "'1238.77641282628"
$ws.Cells.Item(53, 5).Value = "'1290.10644550925"
$ws.Cells.Item(54, 5).Value = "'1110.82169780867"
$ws.Cells.Item(55, 5).Value = "'1198.44325057469"
$ws.Cells.Item(56, 5).Value = "'1239.89697488554"
$ws.Cells.Item(57, 5).Value = "'1276.48678813932"
$ws.Cells.Item(58, 5).Value = "'1322.47951936824"
$ws.Cells.Item(59, 5).Value = "'1384.18500367718"
$ws.Cells.Item(60, 5).Value = "'1458.785339507"
$ws.Cells.Item(61, 5).Value = "'1367.68625725195"
$ws.Cells.Item(62, 5).Value = "'1349.61983373752"

# --- Append new rows for years 2011-2016 (rows 63-68) ---
$ws.Cells.Item(63, 1).Value = 450
$ws.Cells.Item(63, 2).Value = "Madagascar"
$ws.Cells.Item(63, 3).Value = "GDP per Capita"
$ws.Cells.Item(63, 4).Value = 2011
$ws.Cells.Item(63, 5).Value = "'1347"
$ws.Cells.Item(64, 1).Value = 450
$ws.Cells.Item(64, 2).Value = "Madagascar"
$ws.Cells.Item(64, 3).Value = "GDP per Capita"
$ws.Cells.Item(64, 4).Value = 2012
$ws.Cells.Item(64, 5).Value = "'1351"
$ws.Cells.Item(65, 1).Value = 450
$ws.Cells.Item(65, 2).Value = "Madagascar"
$ws.Cells.Item(65, 3).Value = "GDP per Capita"
$ws.Cells.Item(65, 4).Value = 2013
$ws.Cells.Item(65, 5).Value = "'1345"
$ws.Cells.Item(66, 1).Value = 450
$ws.Cells.Item(66, 2).Value = "Madagascar"
$ws.Cells.Item(66, 3).Value = "GDP per Capita"
$ws.Cells.Item(66, 4).Value = 2014
$ws.Cells.Item(66, 5).Value = "'1354"
$ws.Cells.Item(67, 1).Value = 450
$ws.Cells.Item(67, 2).Value = "Madagascar"
$ws.Cells.Item(67, 3).Value = "GDP per Capita"
$ws.Cells.Item(67, 4).Value = 2015
$ws.Cells.Item(67, 5).Value = "'1360"
$ws.Cells.Item(68, 1).Value = 450
$ws.Cells.Item(68, 2).Value = "Madagascar"
$ws.Cells.Item(68, 3).Value = "GDP per Capita"
$ws.Cells.Item(68, 4).Value = 2016
$ws.Cells.Item(68, 5).Value = "'1381"
